$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 86: 四方坪站 (四方坪站 = shared string index 2 already present) ---
$ws.Cells.Item(86, 1).Value = 45943
$ws.Cells.Item(86, 2).Value = "四方坪站"
$ws.Cells.Item(86, 3).Formula = "=15232/126"
$ws.Cells.Item(86, 4).Formula = "=C86/(24*60)"
$ws.Cells.Item(86, 5).Formula = "=8434.09/126"
$ws.Cells.Item(86, 6).Formula = "=2913.27/126"
$ws.Cells.Item(86, 7).Formula = "=8434.09/(15232/60)"
$ws.Cells.Item(86, 8).Formula = "=369/126"

# --- Row 87: 高岭站 (高岭站 = shared string index 3 already present) ---
$ws.Cells.Item(87, 1).Value = 45943
$ws.Cells.Item(87, 2).Value = "高岭站"
$ws.Cells.Item(87, 3).Formula = "=6611/36"
$ws.Cells.Item(87, 4).Formula = "=C87/(24*60)"
$ws.Cells.Item(87, 5).Formula = "=4333.89/36"
$ws.Cells.Item(87, 6).Formula = "=1107.18/36"
$ws.Cells.Item(87, 7).Formula = "=4333.89/(6611/60)"
$ws.Cells.Item(87, 8).Formula = "=169/36"

# --- Update the view/selection to match the authored state ---
$ws.Range("A79").Select()
$excel.ActiveWindow.ScrollRow = 79
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.TabRatio = 601
$ws.Range("J88").Select()

Write-Host "Rows 86-87 added"
